$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-10 (columns A:T), reflecting the "Natmi following
# Dr Hou advice" update: the Sending-cluster combinations now span all three
# clusters (ECs, FAPs, sCs) crossed with all three Target clusters, and the
# associated Ligand/Receptor statistics are recomputed accordingly.
$data = @(
  @("ECs","Angpt1","Itgb1","ECs",2,0.6666666666666666,0.1659303333333333,0.497791,0.01127931880021135,0.01127931880021135,3,1,117.044563,351.133689,0.3245365645427815,0.3245365645427815,19.42124335344433,174.791190180999,0.0036605513738034,0.0036605513738034),
  @("ECs","Angpt1","Itgb1","FAPs",2,0.6666666666666666,0.1659303333333333,0.497791,0.01127931880021135,0.01127931880021135,3,1,101.5800373333333,304.740112,0.281657135515876,0.281657135515876,16.85520945473244,151.696885092592,0.003176900623837897,0.003176900623837896),
  @("ECs","Angpt1","Itgb1","sCs",2,0.6666666666666666,0.1659303333333333,0.497791,0.01127931880021135,0.01127931880021135,3,1,142.0267893333333,426.080368,0.3938062999413425,0.3938062999413425,23.56655249634311,212.098972467088,0.004441866802570056,0.004441866802570056),
  @("FAPs","Angpt1","Itgb1","ECs",3,1,12.95459633333333,38.863789,0.8806046431437031,0.8806046431437031,3,1,117.044563,351.133689,0.3245365645427815,0.3245365645427815,1516.265066676402,13646.38560008762,0.2857884056062794,0.2857884056062794),
  @("FAPs","Angpt1","Itgb1","FAPs",3,1,12.95459633333333,38.863789,0.8806046431437031,0.8806046431437031,3,1,101.5800373333333,304.740112,0.281657135515876,0.281657135515876,1315.928379178263,11843.35541260437,0.2480285813098356,0.2480285813098356),
  @("FAPs","Angpt1","Itgb1","sCs",3,1,12.95459633333333,38.863789,0.8806046431437031,0.8806046431437031,3,1,142.0267893333333,426.080368,0.3938062999413425,0.3938062999413425,1839.899724332706,16559.09751899435,0.3467876562275881,0.3467876562275881),
  @("sCs","Angpt1","Itgb1","ECs",3,1,1.590497666666667,4.771493,0.1081160380560855,0.1081160380560855,3,1,117.044563,351.133689,0.3245365645427815,0.3245365645427815,186.1591043475196,1675.431939127677,0.03508760756269861,0.03508760756269861),
  @("sCs","Angpt1","Itgb1","FAPs",3,1,1.590497666666667,4.771493,0.1081160380560855,0.1081160380560855,3,1,101.5800373333333,304.740112,0.281657135515876,0.281657135515876,161.5628123585795,1454.065311227216,0.03045165358220248,0.03045165358220248),
  @("sCs","Angpt1","Itgb1","sCs",3,1,1.590497666666667,4.771493,0.1081160380560855,0.1081160380560855,3,1,142.0267893333333,426.080368,0.3938062999413425,0.3938062999413425,225.8932770388249,2033.039493349424,0.04257677691118442,0.04257677691118442)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
